# Normalizes the "Recorded By" (column G) values so that the literal
# token "System" (capital S, exact case match) is moved to the end of
# the comma-separated list of recorders, e.g.:
#   "System, dnasr281@gmail.com"            -> "dnasr281@gmail.com, System"
#   "System, admin@admin.com"               -> "admin@admin.com, System"
#   "backup@backdoor.com, System, system"   -> "backup@backdoor.com, system, System"
# Values that do not contain "System" as a standalone item (or already
# have it at the end) are left untouched.
#
# NOTE: this engine's case-sensitive operators (-ceq/-cne/-cmatch/-clike)
# do not actually perform case-sensitive comparisons, so we use a .NET
# Regex object directly (which does honor case) to find the exact-case
# "System" token.

$sysRegex = New-Object System.Text.RegularExpressions.Regex("^System$")

function Transform-RecordedBy($s) {
    $rawParts = $s.Split(",")
    $n = $rawParts.Length
    $trimmed = @()
    $sysIndex = -1
    for ($i = 0; $i -lt $n; $i++) {
        $t = $rawParts[$i].Trim()
        $trimmed += $t
        if ($sysIndex -lt 0 -and $sysRegex.IsMatch($t)) {
            $sysIndex = $i
        }
    }
    # Nothing to do if there's no exact "System" token, or it's already
    # the last element in the list.
    if ($sysIndex -lt 0 -or $sysIndex -eq ($n - 1)) {
        return $null
    }
    $result = @()
    for ($i = 0; $i -lt $n; $i++) {
        if ($i -ne $sysIndex) {
            $result += $trimmed[$i]
        }
    }
    $result += "System"
    return [string]::Join(", ", $result)
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$firstRow = $used.Row
$rowCount = $used.Rows.Count
$lastRow = $firstRow + $rowCount - 1

for ($r = $firstRow; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # Column G = "Recorded By"
    $val = $cell.Value2
    if ($val -ne $null -and $val -is [string] -and $val.Contains(",")) {
        $newVal = Transform-RecordedBy $val
        if ($newVal -ne $null) {
            $cell.Value = $newVal
        }
    }
}
